$d = $word.ActiveDocument

# Locate the paragraph containing the unique marker text that is being
# removed ("//We want two others...") so the edit is robust to any
# paraId/ordering differences.
$rng = $d.Content
$targetText = "//We want two others in the background but still deciding on who they will be so right now we" + [char]0x2019 + "ll just have Taylor Swift and Drake"
$found = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text"
}

$targetPara = $rng.Paragraphs.Item(1)
$targetIndex = $targetPara.Index

# The 3 empty paragraphs immediately preceding the target paragraph are
# removed along with it.
$startPara = $d.Paragraphs.Item($targetIndex - 3)

$delStart = $startPara.Range.Start
$delEnd = $targetPara.Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

$insPoint = $d.Range($delStart, $delStart)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>----------Opening Scene Ends and Transitions to Junk Cave and our mc-----------------</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>[Main Character (yet unnamed)]</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Maybe my parents were righ</w:t></w:r><w:r><w:t>t.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Maybe </w:t></w:r><w:r><w:t>he</w:t></w:r><w:r><w:t xml:space="preserve"> is gone</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>It’s been 4 years since I left home to search for him…</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>My parents thought I was craz</w:t></w:r><w:r><w:t>y</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>They said that he was long gone, and that the world was better off like that</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">That’s when I </w:t></w:r><w:r><w:t>found it…</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>I found a singular Yeezy, a shoe society had left behind after the accident</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Somehow, I knew that this Yeezy was his…</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">That finding it’s match would lead me to the truth… </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>That’s how I ended up here.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Living i</w:t></w:r><w:r><w:t xml:space="preserve">n a Junkyard </w:t></w:r><w:r><w:t>searching</w:t></w:r><w:r><w:t xml:space="preserve"> high and low for a shoe that might not even exist.</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insPoint.InsertXML($xml)
